# Add two new columns (I: "I0", J: "IF") to the right of the existing
# data (A:H), mirroring the header style used by the other header cells
# and filling in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting already used by the other header cells
# (bold font, thin box border, centered / top-aligned).
$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Data ----------------------------------------------------------------
$I = @(6,1,1,1,8,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,6,1,1,1,1)
$J = @(8,5,6,6,9,5,6,4,6,4,5,5,6,6,7,6,6,7,7,6,6,5,7,3,6,6,6,6,6,6,9,9,4,3,3,2)

for ($i = 0; $i -lt $I.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I[$i]
    $ws.Cells.Item($row, 10).Value = $J[$i]
}
